# Remove the two white "filler" rectangles (object 16 / object 17) that
# were left over on the "Using a Lifecycle Hook" slide (slide 16).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

$s.Shapes.Item("object 16").Delete()
$s.Shapes.Item("object 17").Delete()
